$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume/percent-change (E) columns with latest scraped values.
# Column D values must remain plain text (prices use "." as thousands separators in the source
# data, e.g. "30.091.10"), so we force text format, assign the value, then clear the number
# format override so the cell keeps its original (default) style but stores a text value.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "30.091.10"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -0.05%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.878.13"
$cell.ClearFormats()

$ws.Range("E4").Value = "  +0.29%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "319.57"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -3.12%  "

$ws.Range("E6").Value = "  +0.24%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5038"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -3.36%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3957"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -3.02%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.08217"
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "42.12"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -1.93%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.093"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -2.95%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "23.60"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +5.78%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.882.52"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -1.97%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.296"
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.196"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -2.83%  "

$ws.Range("E16").Value = "  +0.23%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "91.65"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -3.83%  "

$ws.Range("E18").Value = "  -2.41%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06464"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -3.41%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "18.07"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -1.92%  "

$ws.Range("E21").Value = "  +0.31%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "30.090.43"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -0.11%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.835"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -3.00%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.16"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -1.90%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.171"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -1.36%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.097.47"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -2.34%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "160.98"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +0.76%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "21.18"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +0.50%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "127.32"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -1.27%  "

$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("E32").Value = "  -2.56%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.929"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -2.15%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.691"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +1.48%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.02423"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -2.61%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.289"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +2.19%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.06349"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -4.00%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.2132"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -3.47%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.173"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -4.77%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "8.497"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -4.68%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.6294"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -3.87%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.212"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -2.91%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "11.29"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -3.04%  "

$ws.Range("E44").Value = "  +0.13%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "13.08"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.87%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.5904"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -4.06%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.092"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +0.62%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "3.626"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -3.56%  "

$ws.Range("E49").Value = "  -3.17%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "122.15"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -1.75%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "77.43"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -2.92%  "
